# Has2-Hmmr LR-pairs sheet: refresh with newly computed TPM-based NATMI output.
#
# The new run:
#   - recomputes every statistic for the existing "FAPs" / "MuSCs" sending
#     clusters (rows 2-13 in the old layout), and
#   - adds a brand-new "ECs" sending cluster block (6 more target-cluster
#     rows), growing the table from 12 data rows to 18 data rows.
#
# Row 1 (the header) is untouched; data rows 2-19 are rewritten in full below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T (Sending cluster, Ligand symbol, Receptor symbol, Target
# cluster, then the 16 numeric NATMI metrics), one array per data row,
# in final row order (rows 2..19).
$data = @(
    @("ECs",   "Has2", "Hmmr", "ECs",              2, 1, 1.459709,           2.919418,           0.02561255265273268, 0.01767271725384302, 2, 1, 3.9645265,          7.929053,           0.1068120917322893,  0.0832154569269431,  5.787055012788501,  23.148220051154,   0.002735730323441772,  0.001470643241419218),
    @("ECs",   "Has2", "Hmmr", "FAPs",             2, 1, 1.459709,           2.919418,           0.02561255265273268, 0.01767271725384302, 3, 1, 18.177566,          54.532698,          0.4897391522195003,  0.5723209797600036,  26.533956688294,    159.203740129764,  0.01254346982232662,   0.01011446685374096),
    @("ECs",   "Has2", "Hmmr", "Inflammatory-Mac", 2, 1, 1.459709,           2.919418,           0.02561255265273268, 0.01767271725384302, 3, 1, 0.5323366666666667, 1.59701,            0.01434219013858555, 0.01676062915292626, 0.7770566233633335, 4.66233974018,     0.0003673401000800259, 0.0002962058600161843),
    @("ECs",   "Has2", "Hmmr", "MuSCs",            2, 1, 1.459709,           2.919418,           0.02561255265273268, 0.01767271725384302, 2, 1, 12.1025545,         24.205109,          0.3260664700939772,  0.2540327584393071,  17.6662077266405,   70.664830906562,   0.008351394633572679,  0.004489449113111679),
    @("ECs",   "Has2", "Hmmr", "Neutrophils",      2, 1, 1.459709,           2.919418,           0.02561255265273268, 0.01767271725384302, 3, 1, 0.7886313333333334, 2.365894,           0.02124726933190069, 0.02483007116369548, 1.151172254948667,  6.907033529692001, 0.0005441968044900988, 0.0004388148270687912),
    @("ECs",   "Has2", "Hmmr", "Resolving-Mac",    2, 1, 1.459709,           2.919418,           0.02561255265273268, 0.01767271725384302, 3, 1, 1.551217333333333, 4.653652,           0.04179282648374708, 0.04884010455712462, 2.264325902422667,  13.585955414536,   0.001070420968821493,  0.0008631373584861936),
    @("FAPs",  "Has2", "Hmmr", "ECs",              3, 1, 51.209624,          153.628872,         0.898541552478366,   0.9299934496816972,  2, 1, 3.9645265,          7.929053,           0.1068120917322893,  0.0832154569269431,  203.021911403036,   1218.131468418216, 0.09597510272859283,   0.07738982985432649),
    @("FAPs",  "Has2", "Hmmr", "FAPs",             3, 1, 51.209624,          153.628872,         0.898541552478366,   0.9299934496816972,  3, 1, 18.177566,          54.532698,          0.4897391522195003,  0.5723209797600036,  930.866320095184,   8377.796880856657, 0.4400509781447486,    0.5322547622922146),
    @("FAPs",  "Has2", "Hmmr", "Inflammatory-Mac", 3, 1, 51.209624,          153.628872,         0.898541552478366,   0.9299934496816972,  3, 1, 0.5323366666666667, 1.59701,            0.01434219013858555, 0.01676062915292626, 27.26076054141333,  245.34684487272,   0.01288705379306457,   0.01558727532476551),
    @("FAPs",  "Has2", "Hmmr", "MuSCs",            3, 1, 51.209624,          153.628872,         0.898541552478366,   0.9299934496816972,  2, 1, 12.1025545,         24.205109,          0.3260664700939772,  0.2540327584393071,  619.7672653845079,  3718.603592307048, 0.2929842722493829,    0.2362488013531285),
    @("FAPs",  "Has2", "Hmmr", "Neutrophils",      3, 1, 51.209624,          153.628872,         0.898541552478366,   0.9299934496816972,  3, 1, 0.7886313333333334, 2.365894,           0.02124726933190069, 0.02483007116369548, 40.38551405461867,  363.469626491568,  0.01909155437141202,   0.02309180353736719),
    @("FAPs",  "Has2", "Hmmr", "Resolving-Mac",    3, 1, 51.209624,          153.628872,         0.898541552478366,   0.9299934496816972,  3, 1, 1.551217333333333, 4.653652,           0.04179282648374708, 0.04884010455712462, 79.43725638228267,  714.9353074405441, 0.03755259119116507,   0.0454209773198951),
    @("MuSCs", "Has2", "Hmmr", "ECs",              2, 1, 4.322604500000001,  8.645209000000001,  0.07584589486890143, 0.05233383306445977, 2, 1, 3.9645265,          7.929053,           0.1068120917322893,  0.0832154569269431,  17.13708008926925,  68.548320357077,   0.008101258680254667,  0.004354983831197382),
    @("MuSCs", "Has2", "Hmmr", "FAPs",             2, 1, 4.322604500000001,  8.645209000000001,  0.07584589486890143, 0.05233383306445977, 3, 1, 18.177566,          54.532698,          0.4897391522195003,  0.5723209797600036,  78.57442859064702,  471.4465715438821, 0.03714470425242513,   0.02995175061404809),
    @("MuSCs", "Has2", "Hmmr", "Inflammatory-Mac", 2, 1, 4.322604500000001,  8.645209000000001,  0.07584589486890143, 0.05233383306445977, 3, 1, 0.5323366666666667, 1.59701,            0.01434219013858555, 0.01676062915292626, 2.301080870848334,  13.80648522509,    0.001087796245440954,  0.0008771479681445605),
    @("MuSCs", "Has2", "Hmmr", "MuSCs",            2, 1, 4.322604500000001,  8.645209000000001,  0.07584589486890143, 0.05233383306445977, 2, 1, 12.1025545,         24.205109,          0.3260664700939772,  0.2540327584393071,  52.31455654319526,  209.258226172781,  0.02473080321102159,   0.01329450797306693),
    @("MuSCs", "Has2", "Hmmr", "Neutrophils",      2, 1, 4.322604500000001,  8.645209000000001,  0.07584589486890143, 0.05233383306445977, 3, 1, 0.7886313333333334, 2.365894,           0.02124726933190069, 0.02483007116369548, 3.408941350307667,  20.453648101846,   0.001611518155998573,  0.001299452799259495),
    @("MuSCs", "Has2", "Hmmr", "Resolving-Mac",    2, 1, 4.322604500000001,  8.645209000000001,  0.07584589486890143, 0.05233383306445977, 3, 1, 1.551217333333333, 4.653652,           0.04179282648374708, 0.04884010455712462, 6.705299025544668,  40.23179415326801, 0.003169814323760521,  0.002555989878743321)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
